$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column: several new values look numeric ("1.00", "0.0000253", ...).
# Force Text format first so Excel keeps the literal string instead of
# auto-converting to a number, then restore the default (unstyled) cell so
# the only thing that changes is the cell text, as in the source workbook.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.700.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.789.66'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.775.17'
$ws.Range("D7").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000253'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.425.04'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.784.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.59'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.693.37'
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '459.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.697'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000154'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.02'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.933.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.78'
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.66'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.10'
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.37'
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.997'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.77'
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '150.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '392.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.90'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.82'
$ws.Range("D51").Style = "Normal"

# Volume(1h) column: percentage text with padding spaces -- already safe
# from Excel's numeric auto-detection, so these can be set directly.
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("E3").Value = '  +1.35%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("E7").Value = '  +1.03%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("E10").Value = '  +0.04%  '
$ws.Range("E11").Value = '  -2.34%  '
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("E13").Value = '  -2.30%  '
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("E15").Value = '  +1.36%  '
$ws.Range("E16").Value = '  +1.28%  '
$ws.Range("E17").Value = '  +4.71%  '
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("E19").Value = '  +0.91%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("E21").Value = '  -5.95%  '
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("E24").Value = '  +5.45%  '
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("E27").Value = '  -2.39%  '
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("E30").Value = '  +1.24%  '
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("E32").Value = '  +3.64%  '
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  -0.42%  '
$ws.Range("E37").Value = '  -0.52%  '
$ws.Range("E38").Value = '  -1.12%  '
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("E44").Value = '  +5.36%  '
$ws.Range("E45").Value = '  +3.33%  '
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("E47").Value = '  +3.99%  '
$ws.Range("E48").Value = '  -1.66%  '
$ws.Range("E49").Value = '  +0.95%  '
$ws.Range("E50").Value = '  +7.41%  '
$ws.Range("E51").Value = '  -4.84%  '
